$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 11, "2021年") below the existing last data row (row 10).
# First copy the formatting of the whole row 10 down to row 11 -- this both carries
# over the header-cell style (bold/border/center alignment) onto A11 and creates the
# same "blank" cells in the unused columns that the rest of the table uses.
$ws.Range("A10:O10").Copy()
$ws.Range("A11:O11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Now fill in the actual values for the new row.
$ws.Cells.Item(11, 1).Value = "2021年"
$ws.Cells.Item(11, 3).Value = 9
$ws.Cells.Item(11, 5).Value = 5
$ws.Cells.Item(11, 7).Value = 9
$ws.Cells.Item(11, 15).Value = 23
